$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.893.56'
$ws.Range("E2").Value = '  +4.19%  '

# Row 3
$ws.Range("D3").Value = '2.346.73'
$ws.Range("E3").Value = '  +2.98%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.95'
$ws.Range("E5").Value = '  +2.94%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.12'
$ws.Range("E6").Value = '  +4.41%  '

# Row 7
$ws.Range("E7").Value = '  +0.36%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.536'
$ws.Range("E8").Value = '  +1.51%  '

# Row 9
$ws.Range("D9").Value = '2.344.03'
$ws.Range("E9").Value = '  +2.13%  '

# Row 10
$ws.Range("E10").Value = '  +7.34%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("E11").Value = '  +0.12%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.27'
$ws.Range("E12").Value = '  +7.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  -0.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.75'
$ws.Range("E14").Value = '  +1.97%  '

# Row 15
$ws.Range("D15").Value = '2.760.12'
$ws.Range("E15").Value = '  +2.86%  '

# Row 16
$ws.Range("D16").Value = '56.799.35'
$ws.Range("E16").Value = '  +3.90%  '

# Row 17
$ws.Range("E17").Value = '  +2.84%  '

# Row 18
$ws.Range("D18").Value = '2.340.04'
$ws.Range("E18").Value = '  +1.76%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.44'
$ws.Range("E19").Value = '  +1.11%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.26'
$ws.Range("E20").Value = '  +3.26%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.21'
$ws.Range("E21").Value = '  +4.68%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("E22").Value = '  +2.55%  '

# Row 23
$ws.Range("E23").Value = '  +0.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.73'
$ws.Range("E24").Value = '  +0.69%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.51%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.160'
$ws.Range("E26").Value = '  +6.84%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.75'
$ws.Range("E27").Value = '  +4.12%  '

# Row 28
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.24'
$ws.Range("E28").Value = '  +11.17%  '

# Row 29
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.13'
$ws.Range("E29").Value = '  -0.30%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0740'
$ws.Range("E30").Value = '  +5.47%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.68'
$ws.Range("E31").Value = '  +3.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.21'
$ws.Range("E32").Value = '  +2.63%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.27'
$ws.Range("E33").Value = '  +1.94%  '

# Row 34
$ws.Range("E34").Value = '  +0.05%  '

# Row 35
$ws.Range("E35").Value = '  +0.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.935'
$ws.Range("E36").Value = '  +2.78%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.25'
$ws.Range("E37").Value = '  +4.01%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.98'
$ws.Range("E38").Value = '  +5.49%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.53'
$ws.Range("E39").Value = '  +8.05%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.53'
$ws.Range("E40").Value = '  +3.06%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.379'
$ws.Range("E41").Value = '  +1.25%  '

# Row 42
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.58'
$ws.Range("E42").Value = '  +5.80%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.94'
$ws.Range("E43").Value = '  +9.00%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '276.60'
$ws.Range("E44").Value = '  +10.61%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.08'
$ws.Range("E45").Value = '  +5.46%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0933'
$ws.Range("E46").Value = '  +3.57%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0504'
$ws.Range("E47").Value = '  +1.73%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.559'
$ws.Range("E48").Value = '  +2.14%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0217'
$ws.Range("E49").Value = '  +5.16%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.379'
$ws.Range("E50").Value = '  +1.42%  '

# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.23'
$ws.Range("E51").Value = '  +6.53%  '
